# This edit re-shuffles the per-row "observation" data (Fecha, Volumen,
# Precio minimo/maximo/promedio, Origen, Precio $/Kg) across the existing
# rows of the sheet, while the row's identity columns (Mercado, Region,
# Categoria, Variedad, Calidad, Unidad, Kg o Unidades, Clasificacion)
# stay put. Row 21 is left untouched.
#
# Because this is a permutation (several rows both donate and receive
# data), we first snapshot every source row's values from the live sheet
# into memory, and only then write them back out to their destination
# rows - avoiding any read-after-write clobbering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps destination row -> source row (source row's D/J/K/L/M/O/P values
# move into the destination row).
$mapping = @{
    2  = 19
    3  = 8
    4  = 18
    5  = 7
    6  = 2
    7  = 26
    8  = 14
    9  = 20
    10 = 24
    11 = 23
    12 = 13
    13 = 15
    14 = 3
    15 = 16
    16 = 4
    17 = 6
    18 = 25
    19 = 5
    20 = 11
    22 = 12
    23 = 22
    24 = 17
    25 = 10
    26 = 9
}

# Snapshot the columns that move: D (Fecha), J (Volumen), K (Precio
# minimo), L (Precio maximo), M (Precio promedio ponderado), O (Origen),
# P (Precio $/Kg).
$snapshot = @{}
foreach ($srcRow in $mapping.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $snapshot[$srcRow] = @{
            D = $ws.Cells.Item($srcRow, 4).Value2
            J = $ws.Cells.Item($srcRow, 10).Value2
            K = $ws.Cells.Item($srcRow, 11).Value2
            L = $ws.Cells.Item($srcRow, 12).Value2
            M = $ws.Cells.Item($srcRow, 13).Value2
            O = $ws.Cells.Item($srcRow, 15).Value2
            P = $ws.Cells.Item($srcRow, 16).Value2
        }
    }
}

# Now write the snapshotted values into the destination rows.
foreach ($dstRow in $mapping.Keys) {
    $srcRow = $mapping[$dstRow]
    $vals = $snapshot[$srcRow]

    $ws.Cells.Item($dstRow, 4).Value = $vals.D
    $ws.Cells.Item($dstRow, 10).Value = $vals.J
    $ws.Cells.Item($dstRow, 11).Value = $vals.K
    $ws.Cells.Item($dstRow, 12).Value = $vals.L
    $ws.Cells.Item($dstRow, 13).Value = $vals.M
    $ws.Cells.Item($dstRow, 15).Value = $vals.O
    $ws.Cells.Item($dstRow, 16).Value = $vals.P
}
